# exemple-valide.xlsx — "Correction des fichiers exemples, champ protection"
#
# - acces (K2) switches from a numeric flag to the free-text value "LIBRE"
# - date_maj (S2) switches from a real date to the literal text "2020-09-17"
#   (format becomes plain Text "@", font/wrap revert to the sheet default)
# - the "protection" column (I) is widened, and the two trailing columns
#   (S/T) get distinct widths instead of sharing one
# - row 2 no longer needs the extra height that wrapping used to require
# - the remembered selection goes back to the top-left cell, A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- acces: now a free-text "LIBRE" instead of a 1/0 flag -----------------
$ws.Range("K2").Value = "LIBRE"

# --- date_maj: store the date as literal text, not a real date value ------
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Font.Name = "Arial"
$ws.Range("S2").WrapText = $false
$ws.Range("S2").Value = "2020-09-17"

# --- row 2 no longer needs the taller, wrapped-text row height ------------
$ws.Rows("2").RowHeight = 12.8

# --- column widths: "protection" widened, S/T split apart -----------------
$ws.Columns("I").ColumnWidth = 27.5
$ws.Columns("S").ColumnWidth = 9.67
$ws.Columns("T").ColumnWidth = 12.17

# --- restore the saved selection to A1 -------------------------------------
[void]$ws.Range("A1").Select()
